$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition the two charts first (while column widths are still the
#     originals) so the twoCellAnchor offsets land exactly where the
#     author dragged them to. ---
$co1 = $ws.ChartObjects().Item(1)
$co1.Left = 261.46015625
$co1.Top = 10.2
$co1.Width = 754.8875
$co1.Height = 252.60000000000008

$co2 = $ws.ChartObjects().Item(2)
$co2.Left = 259.66015625
$co2.Top = 272.40000000000003
$co2.Width = 760.8875
$co2.Height = 257.69999999999965

# --- Column B (Number Rolled instances) updates for the new session's data ---
$ws.Range("B7").Value = 1
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B11").Value = 3
$ws.Range("B12").Value = 1
$ws.Range("B13").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B19").Value = 2
$ws.Range("B21").Value = 2

# --- New "Portent" mini-table headers (Q1:T1), matching the formatting of
#     the existing Q1/R1/S1 header cells ---
$ws.Range("Q1").Value = "Portent:"
$ws.Range("Q1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("R1").Value = "Number Rolled"
$ws.Range("S1").Value = "Instances Rolled"
$ws.Range("T1").Value = "%"

# --- Column R (Number Rolled, 1-20) ---
$ws.Range("R2").Value = 1
$ws.Range("R3").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("R5").Value = 4
$ws.Range("R6").Value = 5
$ws.Range("R7").Value = 6
$ws.Range("R8").Value = 7
$ws.Range("R9").Value = 8
$ws.Range("R10").Value = 9
$ws.Range("R11").Value = 10
$ws.Range("R12").Value = 11
$ws.Range("R13").Value = 12
$ws.Range("R14").Value = 13
$ws.Range("R15").Value = 14
$ws.Range("R16").Value = 15
$ws.Range("R17").Value = 16
$ws.Range("R18").Value = 17
$ws.Range("R19").Value = 18
$ws.Range("R20").Value = 19
$ws.Range("R21").Value = 20

# --- Column S (Instances Rolled for Portent) ---
$ws.Range("S8").Value = 1
$ws.Range("S21").Value = 1
$ws.Range("S22").Formula = "=SUM(S2:S21)"

# --- Column T (%) ---
$ws.Range("T2").Formula = "=S2/S22"
$ws.Range("T3").Formula = "=S3/S22"
$ws.Range("T4").Formula = "=S4/S22"
$ws.Range("T5").Formula = "=S5/S22"
$ws.Range("T6").Formula = "=S6/S22"
$ws.Range("T7").Formula = "=S7/S22"
$ws.Range("T8").Formula = "=S8/S22"
$ws.Range("T9").Formula = "=S9/S22"
$ws.Range("T10").Formula = "=S10/S22"
$ws.Range("T11").Formula = "=S11/S22"
$ws.Range("T12").Formula = "=S12/S22"
$ws.Range("T13").Formula = "=S13/S22"
$ws.Range("T14").Formula = "=S14/S22"
$ws.Range("T15").Formula = "=S15/S22"
$ws.Range("T16").Formula = "=S16/S22"
$ws.Range("T17").Formula = "=S17/S22"
$ws.Range("T18").Formula = "=S18/S22"
$ws.Range("T19").Formula = "=S19/S22"
$ws.Range("T20").Formula = "=S20/S22"
$ws.Range("T21").Formula = "=S21/S22"

# --- Column widths (set after the charts are positioned) ---
$ws.Columns.Item(16).AutoFit()
$ws.Columns.Item(19).ColumnWidth = 14.6
$ws.Columns.Item(20).ColumnWidth = 12.6

# --- Selection ---
$ws.Range("B13").Select()
